{"js": "const pairs = [\n  [\"14+11=25\", \"71+18=89\"],\n  [\"50-32=18\", \"22+38=60\"],\n  [\"96-77=19\", \"38-38=0\"],\n  [\"24+74=98\", \"13+21=34\"],\n  [\"72-44=28\", \"5+15=20\"],\n  [\"54-30=24\", \"21+42=63\"],\n  [\"74-25=49\", \"19+24=43\"],\n  [\"25+18=43\", \"43+55=98\"],\n  [\"91-11=80\", \"50-16=34\"],\n  [\"95-51=44\", \"61-39=22\"],\n  [\"8+38=46\", \"50+16=66\"],\n  [\"50+49=99\", \"14+17=31\"],\n  [\"45+22=67\", \"23+42=65\"],\n  [\"54+40=94\", \"12+53=65\"],\n  [\"8+24=32\", \"10+47=57\"],\n  [\"22+16=38\", \"39-4=35\"],\n  [\"47+43=90\", \"88-68=20\"],\n  [\"2+17=19\", \"38-29=9\"],\n  [\"87-33=54\", \"52+31=83\"],\n  [\"30+0=30\", \"16+46=62\"],\n  [\"58+2=60\", \"14+42=56\"],\n  [\"37+43=80\", \"11-0=11\"],\n  [\"33-19=14\", \"78+5=83\"],\n  [\"27+56=83\", \"36-35=1\"],\n  [\"67-49=18\", \"63+33=96\"],\n  [\"63-1=62\", \"88-0=88\"],\n  [\"68+14=82\", \"52+32=84\"],\n  [\"33+39=72\", \"7-7=0\"],\n  [\"66-21=45\", \"2+60=62\"],\n  [\"87-69=18\", \"29+59=88\"],\n  [\"47+44=91\", \"97-72=25\"],\n  [\"77-28=49\", \"53+41=94\"],\n  [\"17+66=83\", \"94-78=16\"],\n  [\"39+5=44\", \"72+7=79\"],\n  [\"12+35=47\", \"9+4=13\"],\n  [\"70-34=36\", \"81-49=32\"],\n  [\"86-6=80\", \"57-31=26\"],\n  [\"25-11=14\", \"95-52=43\"],\n  [\"34+18=52\", \"14+47=61\"],\n  [\"93-90=3\", \"1+14=15\"],\n  [\"43-14=29\", \"77-27=50\"],\n  [\"63-15=48\", \"18-12=6\"],\n  [\"24+40=64\", \"96-95=1\"],\n  [\"46-3=43\", \"2+86=88\"],\n  [\"82-61=21\", \"68+1=69\"],\n  [\"57-32=25\", \"44-33=11\"],\n  [\"80+1=81\", \"20+3=23\"],\n  [\"72-7=65\", \"21+43=64\"],\n  [\"57-1=56\", \"43+8=51\"],\n  [\"66+24=90\", \"96-12=84\"],\n  [\"2+16=18\", \"72-43=29\"],\n  [\"23-8=15\", \"21+63=84\"],\n  [\"51+19=70\", \"59-20=39\"],\n  [\"37+20=57\", \"72-0=72\"],\n  [\"81-48=33\", \"31+23=54\"],\n  [\"99-19=80\", \"77-60=17\"],\n  [\"38+18=56\", \"49+1=50\"],\n  [\"70+9=79\", \"70-12=58\"],\n  [\"32+6=38\", \"72-39=33\"],\n  [\"78-62=16\", \"34-21=13\"],\n  [\"48+13=61\", \"59+15=74\"],\n  [\"20+26=46\", \"36+0=36\"],\n  [\"68+19=87\", \"12+8=20\"],\n  [\"83+15=98\", \"7+79=86\"],\n  [\"41+44=85\", \"50-43=7\"],\n  [\"43-27=16\", \"49+25=74\"],\n  [\"28-22=6\", \"81+8=89\"],\n  [\"46+5=51\", \"22+23=45\"],\n  [\"40+53=93\", \"59+38=97\"],\n  [\"77-51=26\", \"20+60=80\"],\n  [\"64-53=11\", \"6+84=90\"],\n  [\"11+81=92\", \"56+42=98\"],\n  [\"49-24=25\", \"49-10=39\"],\n  [\"16+52=68\", \"55+7=62\"],\n  [\"20+75=95\", \"49-36=13\"],\n  [\"65-10=55\", \"50-15=35\"],\n  [\"50+15=65\", \"44-22=22\"],\n  [\"43+47=90\", \"34+12=46\"],\n  [\"35-3=32\", \"86-59=27\"],\n  [\"40-31=9\", \"72-26=46\"],\n  [\"74-18=56\", \"20+39=59\"],\n  [\"13+20=33\", \"78-12=66\"],\n  [\"41+5=46\", \"18-2=16\"],\n  [\"37+25=62\", \"93-4=89\"],\n  [\"27+9=36\", \"28+14=42\"],\n  [\"32+14=46\", \"67-54=13\"],\n  [\"90-58=32\", \"59+5=64\"],\n  [\"63+31=94\", \"74-66=8\"],\n  [\"77-64=13\", \"70+18=88\"],\n  [\"70-17=53\", \"87-79=8\"],\n  [\"41+52=93\", \"51-16=35\"],\n  [\"92+6=98\", \"16+28=44\"],\n  [\"94-43=51\", \"96-68=28\"],\n  [\"26+19=45\", \"84-55=29\"],\n  [\"83-58=25\", \"5+70=75\"],\n  [\"39+50=89\", \"19+65=84\"],\n  [\"16+6=22\", \"1+86=87\"],\n  [\"82-28=54\", \"21+74=95\"],\n  [\"62+31=93\", \"1+30=31\"],\n  [\"2+80=82\", \"76+14=90\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    console.log(\"WARNING: not found -> \" + oldText);\n    continue;\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n    @(\"14+11=25\", \"71+18=89\"),\n    @(\"50-32=18\", \"22+38=60\"),\n    @(\"96-77=19\", \"38-38=0\"),\n    @(\"24+74=98\", \"13+21=34\"),\n    @(\"72-44=28\", \"5+15=20\"),\n    @(\"54-30=24\", \"21+42=63\"),\n    @(\"74-25=49\", \"19+24=43\"),\n    @(\"25+18=43\", \"43+55=98\"),\n    @(\"91-11=80\", \"50-16=34\"),\n    @(\"95-51=44\", \"61-39=22\"),\n    @(\"8+38=46\", \"50+16=66\"),\n    @(\"50+49=99\", \"14+17=31\"),\n    @(\"45+22=67\", \"23+42=65\"),\n    @(\"54+40=94\", \"12+53=65\"),\n    @(\"8+24=32\", \"10+47=57\"),\n    @(\"22+16=38\", \"39-4=35\"),\n    @(\"47+43=90\", \"88-68=20\"),\n    @(\"2+17=19\", \"38-29=9\"),\n    @(\"87-33=54\", \"52+31=83\"),\n    @(\"30+0=30\", \"16+46=62\"),\n    @(\"58+2=60\", \"14+42=56\"),\n    @(\"37+43=80\", \"11-0=11\"),\n    @(\"33-19=14\", \"78+5=83\"),\n    @(\"27+56=83\", \"36-35=1\"),\n    @(\"67-49=18\", \"63+33=96\"),\n    @(\"63-1=62\", \"88-0=88\"),\n    @(\"68+14=82\", \"52+32=84\"),\n    @(\"33+39=72\", \"7-7=0\"),\n    @(\"66-21=45\", \"2+60=62\"),\n    @(\"87-69=18\", \"29+59=88\"),\n    @(\"47+44=91\", \"97-72=25\"),\n    @(\"77-28=49\", \"53+41=94\"),\n    @(\"17+66=83\", \"94-78=16\"),\n    @(\"39+5=44\", \"72+7=79\"),\n    @(\"12+35=47\", \"9+4=13\"),\n    @(\"70-34=36\", \"81-49=32\"),\n    @(\"86-6=80\", \"57-31=26\"),\n    @(\"25-11=14\", \"95-52=43\"),\n    @(\"34+18=52\", \"14+47=61\"),\n    @(\"93-90=3\", \"1+14=15\"),\n    @(\"43-14=29\", \"77-27=50\"),\n    @(\"63-15=48\", \"18-12=6\"),\n    @(\"24+40=64\", \"96-95=1\"),\n    @(\"46-3=43\", \"2+86=88\"),\n    @(\"82-61=21\", \"68+1=69\"),\n    @(\"57-32=25\", \"44-33=11\"),\n    @(\"80+1=81\", \"20+3=23\"),\n    @(\"72-7=65\", \"21+43=64\"),\n    @(\"57-1=56\", \"43+8=51\"),\n    @(\"66+24=90\", \"96-12=84\"),\n    @(\"2+16=18\", \"72-43=29\"),\n    @(\"23-8=15\", \"21+63=84\"),\n    @(\"51+19=70\", \"59-20=39\"),\n    @(\"37+20=57\", \"72-0=72\"),\n    @(\"81-48=33\", \"31+23=54\"),\n    @(\"99-19=80\", \"77-60=17\"),\n    @(\"38+18=56\", \"49+1=50\"),\n    @(\"70+9=79\", \"70-12=58\"),\n    @(\"32+6=38\", \"72-39=33\"),\n    @(\"78-62=16\", \"34-21=13\"),\n    @(\"48+13=61\", \"59+15=74\"),\n    @(\"20+26=46\", \"36+0=36\"),\n    @(\"68+19=87\", \"12+8=20\"),\n    @(\"83+15=98\", \"7+79=86\"),\n    @(\"41+44=85\", \"50-43=7\"),\n    @(\"43-27=16\", \"49+25=74\"),\n    @(\"28-22=6\", \"81+8=89\"),\n    @(\"46+5=51\", \"22+23=45\"),\n    @(\"40+53=93\", \"59+38=97\"),\n    @(\"77-51=26\", \"20+60=80\"),\n    @(\"64-53=11\", \"6+84=90\"),\n    @(\"11+81=92\", \"56+42=98\"),\n    @(\"49-24=25\", \"49-10=39\"),\n    @(\"16+52=68\", \"55+7=62\"),\n    @(\"20+75=95\", \"49-36=13\"),\n    @(\"65-10=55\", \"50-15=35\"),\n    @(\"50+15=65\", \"44-22=22\"),\n    @(\"43+47=90\", \"34+12=46\"),\n    @(\"35-3=32\", \"86-59=27\"),\n    @(\"40-31=9\", \"72-26=46\"),\n    @(\"74-18=56\", \"20+39=59\"),\n    @(\"13+20=33\", \"78-12=66\"),\n    @(\"41+5=46\", \"18-2=16\"),\n    @(\"37+25=62\", \"93-4=89\"),\n    @(\"27+9=36\", \"28+14=42\"),\n    @(\"32+14=46\", \"67-54=13\"),\n    @(\"90-58=32\", \"59+5=64\"),\n    @(\"63+31=94\", \"74-66=8\"),\n    @(\"77-64=13\", \"70+18=88\"),\n    @(\"70-17=53\", \"87-79=8\"),\n    @(\"41+52=93\", \"51-16=35\"),\n    @(\"92+6=98\", \"16+28=44\"),\n    @(\"94-43=51\", \"96-68=28\"),\n    @(\"26+19=45\", \"84-55=29\"),\n    @(\"83-58=25\", \"5+70=75\"),\n    @(\"39+50=89\", \"19+65=84\"),\n    @(\"16+6=22\", \"1+86=87\"),\n    @(\"82-28=54\", \"21+74=95\"),\n    @(\"62+31=93\", \"1+30=31\"),\n    @(\"2+80=82\", \"76+14=90\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $found = $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n    if (-not $found) {\n        Write-Host \"WARNING: not found -> $old\"\n    }\n}"}
